$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new config column: C1 header "RetryLimit", C2 value 3
$ws.Range("C1").Value = "RetryLimit"
$ws.Range("C2").Value = 3

# Update the active selection to mirror the recorded state in the saved file
$ws.Range("C8").Select()
